# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new (blank) column inserted
# immediately before column N ("Late"), shifting the existing N/O/P
# columns one place to the right (N->O, O->P, P->Q). The new column
# inherits the width/format of its left neighbour (column M), which is
# the standard Excel behaviour when inserting a column.
#
# The workbook's active sheet is also switched from "Transactions" to
# "Repayment schedule", with the selection there moved to cell R8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N - this shifts N,O,P -> O,P,Q
$ws.Columns("N").Insert()

# The newly inserted column picks up column M's width (10.71 chars wide)
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet and select cell R8 on it
# (this also clears the tabSelected flag that used to be on "Transactions")
$ws.Activate()
$ws.Range("R8").Select()
